$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
Set-CellText $ws.Range('D2') '61.482.84'
Set-CellText $ws.Range('E2') '  +1.04%  '

# Row 3
Set-CellText $ws.Range('D3') '3.445.20'
Set-CellText $ws.Range('E3') '  +1.81%  '

# Row 4
Set-CellText $ws.Range('D4') '1.00'

# Row 5
Set-CellText $ws.Range('D5') '579.96'
Set-CellText $ws.Range('E5') '  +0.89%  '

# Row 6
Set-CellText $ws.Range('D6') '149.26'
Set-CellText $ws.Range('E6') '  +9.22%  '

# Row 7
Set-CellText $ws.Range('D7') '3.447.02'
Set-CellText $ws.Range('E7') '  +1.89%  '

# Row 8
Set-CellText $ws.Range('E8') '  +0.12%  '

# Row 9
Set-CellText $ws.Range('E9') '  +0.84%  '

# Row 10
Set-CellText $ws.Range('D10') '7.82'
Set-CellText $ws.Range('E10') '  +3.11%  '

# Row 11
Set-CellText $ws.Range('E11') '  +1.44%  '

# Row 12
Set-CellText $ws.Range('D12') '0.393'
Set-CellText $ws.Range('E12') '  +1.07%  '

# Row 13
Set-CellText $ws.Range('D13') '4.033.77'
Set-CellText $ws.Range('E13') '  +1.86%  '

# Row 14
Set-CellText $ws.Range('D14') '28.00'
Set-CellText $ws.Range('E14') '  +6.71%  '

# Row 15
Set-CellText $ws.Range('E15') '  -0.41%  '

# Row 16
Set-CellText $ws.Range('E16') '  +0.80%  '

# Row 17
Set-CellText $ws.Range('D17') '3.441.11'
Set-CellText $ws.Range('E17') '  +1.75%  '

# Row 18
Set-CellText $ws.Range('D18') '61.593.57'
Set-CellText $ws.Range('E18') '  +1.05%  '

# Row 19
Set-CellText $ws.Range('D19') '6.34'
Set-CellText $ws.Range('E19') '  +8.60%  '

# Row 20
Set-CellText $ws.Range('D20') '14.31'
Set-CellText $ws.Range('E20') '  +1.82%  '

# Row 21
Set-CellText $ws.Range('D21') '9.49'
Set-CellText $ws.Range('E21') '  -0.06%  '

# Row 22
Set-CellText $ws.Range('D22') '388.87'
Set-CellText $ws.Range('E22') '  +2.87%  '

# Row 23
Set-CellText $ws.Range('E23') '  +2.50%  '

# Row 24
Set-CellText $ws.Range('D24') '3.588.63'
Set-CellText $ws.Range('E24') '  +1.80%  '

# Row 25
Set-CellText $ws.Range('D25') '72.85'
Set-CellText $ws.Range('E25') '  +2.17%  '

# Row 26
Set-CellText $ws.Range('B26') 'LEO'
Set-CellText $ws.Range('C26') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-CellText $ws.Range('D26') '5.78'
Set-CellText $ws.Range('E26') '  +0.83%  '

# Row 27
Set-CellText $ws.Range('B27') 'Dai'
Set-CellText $ws.Range('C27') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText $ws.Range('D27') '1.00'
Set-CellText $ws.Range('E27') '  -0.08%  '

# Row 28
Set-CellText $ws.Range('E28') '  -1.52%  '

# Row 29
Set-CellText $ws.Range('D29') '0.181'
Set-CellText $ws.Range('E29') '  +6.21%  '

# Row 30
Set-CellText $ws.Range('D30') '7.80'
Set-CellText $ws.Range('E30') '  +3.40%  '

# Row 31
Set-CellText $ws.Range('D31') '1.00'
Set-CellText $ws.Range('E31') '  +0.04%  '

# Row 32
Set-CellText $ws.Range('E32') '  -13.54%  '

# Row 33
Set-CellText $ws.Range('E33') '  +1.19%  '

# Row 34
Set-CellText $ws.Range('E34') '  +0.73%  '

# Row 36
Set-CellText $ws.Range('E36') '  +1.02%  '

# Row 37
Set-CellText $ws.Range('D37') '5.30'
Set-CellText $ws.Range('E37') '  +1.59%  '

# Row 38
Set-CellText $ws.Range('D38') '7.08'
Set-CellText $ws.Range('E38') '  +3.24%  '

# Row 39
Set-CellText $ws.Range('E39') '  +1.85%  '

# Row 40
Set-CellText $ws.Range('D40') '165.84'
Set-CellText $ws.Range('E40') '  +0.72%  '

# Row 41
Set-CellText $ws.Range('E41') '  +4.74%  '

# Row 42
Set-CellText $ws.Range('D42') '26.50'
Set-CellText $ws.Range('E42') '  +9.77%  '

# Row 43
Set-CellText $ws.Range('E43') '  +2.82%  '

# Row 44
Set-CellText $ws.Range('B44') 'FirstDigitalUSD'
Set-CellText $ws.Range('C44') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText $ws.Range('D44') '1.00'
Set-CellText $ws.Range('E44') '  +0.01%  '

# Row 45
Set-CellText $ws.Range('B45') 'Filecoin'
Set-CellText $ws.Range('C45') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws.Range('D45') '4.50'
Set-CellText $ws.Range('E45') '  +1.79%  '

# Row 46
Set-CellText $ws.Range('D46') '42.27'
Set-CellText $ws.Range('E46') '  +1.63%  '

# Row 47
Set-CellText $ws.Range('E47') '  +0.45%  '

# Row 48
Set-CellText $ws.Range('D48') '2.612.69'
Set-CellText $ws.Range('E48') '  +6.65%  '

# Row 49
Set-CellText $ws.Range('E49') '  -3.42%  '

# Row 50
Set-CellText $ws.Range('D50') '7.07'
Set-CellText $ws.Range('E50') '  +3.72%  '

# Row 51
Set-CellText $ws.Range('D51') '23.19'
Set-CellText $ws.Range('E51') '  -1.22%  '
